$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=19.21311956513562; "C"=5.743344418853512; "D"=5.216813384972417; "F"=49.51549123950535; "G"=3.763195713753014; "I"=35.66441405098071; "J"=10.23034262778265; "K"=16.88250202672441; "L"=12.33274353564821; "N"=24.21248032906135 }
  3 = @{ "B"=19.0876989864023; "C"=5.602500359000041; "D"=5.215097600140226; "F"=49.53268130166698; "G"=3.766032797214222; "I"=35.71591516872429; "J"=10.24814399230999; "K"=16.79865698260188; "L"=12.34064454997506; "N"=24.26824382588179 }
  4 = @{ "B"=19.014841138193; "C"=5.516057409021329; "D"=5.214660658222501; "F"=49.55237377341693; "G"=3.767867048613937; "I"=35.75310255417224; "J"=10.25981989091935; "K"=16.7508054396633; "L"=12.34731203746599; "N"=24.30441262114261 }
  5 = @{ "B"=18.98622182204823; "C"=5.480898985241092; "D"=5.214638478663403; "F"=49.5626957063345; "G"=3.768637801805388; "I"=35.76965479075451; "J"=10.26476588546568; "K"=16.73223418934421; "L"=12.35048637408031; "N"=24.31963768153414 }
  6 = @{ "B"=18.98153503888935; "C"=5.475066800050466; "D"=5.214644234089754; "F"=49.56454835753428; "G"=3.768767193205446; "I"=35.77248767185788; "J"=10.26559852993397; "K"=16.72920697485916; "L"=12.35104110360973; "N"=24.32219516483879 }
  7 = @{ "B"=19.01445079785164; "C"=5.515582896872914; "D"=5.214659726869174; "F"=49.55250367939639; "G"=3.767877348893629; "I"=35.75332012475993; "J"=10.25988583266327; "K"=16.75055120119878; "L"=12.34735299548543; "N"=24.30461598297452 }
  8 = @{ "B"=19.16902980938279; "C"=5.694808621021211; "D"=5.216094331460195; "F"=49.51952169107236; "G"=3.764154836761567; "I"=35.68101572594085; "J"=10.23632601210968; "K"=16.85284887501571; "L"=12.33509122044077; "N"=24.23130749171478 }
  9 = @{ "B"=19.50382333826286; "C"=6.044076350030372; "D"=5.223761562249408; "F"=49.52734715175666; "G"=3.757583543004792; "I"=35.58344453169283; "J"=10.19602323498358; "K"=17.08149445033021; "L"=12.32542859332628; "N"=24.10282840827457 }
  10 = @{ "B"=19.76729196338151; "C"=6.296208850349151; "D"=5.232300799579272; "F"=49.57725144989681; "G"=3.753194737436422; "I"=35.53878619613234; "J"=10.16998152089232; "K"=17.26547685496569; "L"=12.32705490869864; "N"=24.01770204541674 }
  11 = @{ "B"=19.89053798049155; "C"=6.40933523603451; "D"=5.236804376974136; "F"=49.60951343985433; "G"=3.751292440123347; "I"=35.52434935245184; "J"=10.15890373053978; "K"=17.35239254906897; "L"=12.3296780343535; "N"=23.98097838131507 }
  12 = @{ "B"=19.93765981290802; "C"=6.451901757990345; "D"=5.238597663324626; "F"=49.62310069456789; "G"=3.750585553028538; "I"=35.51972817901041; "J"=10.15481896467974; "K"=17.38574547414563; "L"=12.33094096254069; "N"=23.96735912130069 }
  13 = @{ "B"=19.92749180516855; "C"=6.442747091352227; "D"=5.238207557162715; "F"=49.62011357095768; "G"=3.750737195817867; "I"=35.5206858126049; "J"=10.15569379914658; "K"=17.37854314320955; "L"=12.33065699809746; "N"=23.97027950747362 }
  14 = @{ "B"=19.8944058909196; "C"=6.41284291275672; "D"=5.236950157164118; "F"=49.61060383671354; "G"=3.751234014522344; "I"=35.52395221421427; "J"=10.15856546897276; "K"=17.3551278358409; "L"=12.32977654025697; "N"=23.97985216364415 }
  15 = @{ "B"=19.87419747797615; "C"=6.394489007102895; "D"=5.236191373312414; "F"=49.60495716026925; "G"=3.75154008256409; "I"=35.52606312879536; "J"=10.16033878338004; "K"=17.34084185723233; "L"=12.32927230720427; "N"=23.98575307835415 }
  16 = @{ "B"=19.7593027035353; "C"=6.288780152169755; "D"=5.232018840617179; "F"=49.57533507557753; "G"=3.753320946167791; "I"=35.53984801036398; "J"=10.17072091587265; "K"=17.25985954922452; "L"=12.3269212776679; "N"=24.02014223928099 }
  17 = @{ "B"=19.68966164090682; "C"=6.223495093966998; "D"=5.229616882244525; "F"=49.55960867889889; "G"=3.754437521033323; "I"=35.54981059425456; "J"=10.17728662768281; "K"=17.21098818040795; "L"=12.32596053634409; "N"=24.04175094535525 }
  18 = @{ "B"=19.6499275397629; "C"=6.185799628603853; "D"=5.22829367687464; "F"=49.5514634575518; "G"=3.75508861553215; "I"=35.55609410806947; "J"=10.18113542959235; "K"=17.18318330645578; "L"=12.3255852894118; "N"=24.05436805031486 }
  19 = @{ "B"=19.63653059413132; "C"=6.173013114611319; "D"=5.227855715690539; "F"=49.54886035157248; "G"=3.755310590639832; "I"=35.55831661008578; "J"=10.18245101066707; "K"=17.17382207567674; "L"=12.32548873141053; "N"=24.05867235025168 }
  20 = @{ "B"=19.69704202993908; "C"=6.230460176226199; "D"=5.229866545199473; "F"=49.56118964350797; "G"=3.754317742234289; "I"=35.54869279178098; "J"=10.17658020840365; "K"=17.21615925816783; "L"=12.32604446261046; "N"=24.03943117203595 }
  21 = @{ "B"=19.90411207500706; "C"=6.421634227570156; "D"=5.23731711029376; "F"=49.61335992574706; "G"=3.751087721790319; "I"=35.52296983747772; "J"=10.15771900348676; "K"=17.36199373214617; "L"=12.33002784526374; "N"=23.97703265207025 }
  22 = @{ "B"=20.04205703282862; "C"=6.544973211903722; "D"=5.242698153722809; "F"=49.65544084828055; "G"=3.749055205994594; "I"=35.51108806328282; "J"=10.14603399392572; "K"=17.45985734827787; "L"=12.33420204456092; "N"=23.93792537917184 }
  23 = @{ "B"=19.96820643128589; "C"=6.479306037102769; "D"=5.239779759303835; "F"=49.6322525931788; "G"=3.750132840186039; "I"=35.5169784527543; "J"=10.15221189868339; "K"=17.40740006556226; "L"=12.33183089707965; "N"=23.9586446804277 }
  24 = @{ "B"=19.69370440795515; "C"=6.227311764256681; "D"=5.229753492745002; "F"=49.56047209768375; "G"=3.754371865687696; "I"=35.54919641864653; "J"=10.17689934975459; "K"=17.21382050279096; "L"=12.32600596781904; "N"=24.04047933723704 }
  25 = @{ "B"=19.41004406957971; "C"=5.950160196683836; "D"=5.221173007333258; "F"=49.51746904714719; "G"=3.759283776923893; "I"=35.60509804860688; "J"=10.20629757205462; "K"=17.01674671413221; "L"=12.32650675292799; "N"=24.13595465527302 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = $rowData[$col]
  }
}